$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural column changes -------------------------------------------
# Drop the "Input_Feature" column (old column C); everything from old D
# onward shifts left by one.
$ws.Columns("C").Delete()

# Drop the old per-class precision/recall/F-score columns
# (S.PRC, S.RCL, R.PRC, R.RCL, S.FSc, R.FSc) -- after the first delete these
# live at M:R, immediately before the trailing "Accuracy" column.
$ws.Columns("M:R").Delete()

# Make room for the three new ROC/AUC-related metric columns ahead of the
# (now shifted-left) "Accuracy" column.
$ws.Columns("M:O").Insert()

# --- New header row --------------------------------------------------------
$ws.Range("M1").Value = "Sensitivity (TPR)"
$ws.Range("N1").Value = "Specificty(TNR)"
$ws.Range("O1").Value = "1-Specificity(FPR)"

# Size the three new columns to fit their (longer) header text, same as
# Excel auto-fitting a freshly-typed header.
$ws.Columns("M").ColumnWidth = 14.67
$ws.Columns("N").ColumnWidth = 13.83
$ws.Columns("O").ColumnWidth = 15.83

# --- Row 2: fix up B/C values (Input_Feature column removal shifted the
# Antibiotic value into C; "none" belongs in B/max_features) and rewrite the
# summary-stat formulas to use the new column layout -----------------------
$ws.Range("B2").Value = "none"
$ws.Range("C2").Value = "AMX"

$ws.Range("M2").Formula = "=L2/(L2+K2)"
$ws.Range("N2").Formula = "=I2/(I2+J2)"
$ws.Range("O2").Formula = "=1-N2"
$ws.Range("P2").Formula = "=((I2+L2)/SUM(I2:L2))"

# --- Row 3: new data row for the second GB run -----------------------------
$ws.Range("A3").Value = 300
$ws.Range("B3").Value = "none"
$ws.Range("C3").Value = "AMX"
$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 0.8
$ws.Range("I3").Value = 75
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 16
$ws.Range("L3").Value = 119
$ws.Range("M3").Formula = "=L3/(L3+K3)"
$ws.Range("N3").Formula = "=I3/(I3+J3)"
$ws.Range("O3").Formula = "=1-N3"

# --- Selection, matching the saved state in the target file ---------------
$ws.Range("P2").Select()
